$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price/Volume columns so numeric-looking strings
# (e.g. "1.003") are written as text, matching the original inline-string cells,
# then restore the default "Normal" style so no stray number-format style is left behind.
$priceVolRange = $ws.Range("D2:E51")
$priceVolRange.NumberFormat = "@"

$ws.Range('D2').Value = '28.952.64'
$ws.Range('E2').Value = '  -1.54%  '
$ws.Range('D3').Value = '1.911.33'
$ws.Range('E3').Value = '  -1.86%  '
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '324.64'
$ws.Range('E5').Value = '  -0.21%  '
$ws.Range('E6').Value = '  +0.06%  '
$ws.Range('E7').Value = '  -0.85%  '
$ws.Range('D8').Value = '0.3821'
$ws.Range('E8').Value = '  -1.25%  '
$ws.Range('D9').Value = '0.07717'
$ws.Range('E9').Value = '  -1.36%  '
$ws.Range('E10').Value = '  +0.13%  '
$ws.Range('D11').Value = '22.07'
$ws.Range('E11').Value = '  -2.63%  '
$ws.Range('D12').Value = '1.915.47'
$ws.Range('E12').Value = '  -0.38%  '
$ws.Range('D13').Value = '6.940'
$ws.Range('E13').Value = '  -2.00%  '
$ws.Range('D14').Value = '5.663'
$ws.Range('E14').Value = '  -1.46%  '
$ws.Range('D15').Value = '0.07036'
$ws.Range('E15').Value = '  -0.08%  '
$ws.Range('D16').Value = '1.004'
$ws.Range('E16').Value = '  +0.04%  '
$ws.Range('E17').Value = '  -3.36%  '
$ws.Range('D18').Value = '0.000009456'
$ws.Range('E18').Value = '  -3.68%  '
$ws.Range('D19').Value = '16.69'
$ws.Range('E19').Value = '  -1.92%  '
$ws.Range('D20').Value = '1.002'
$ws.Range('E20').Value = '  -0.07%  '
$ws.Range('D21').Value = '28.931.26'
$ws.Range('E21').Value = '  -1.71%  '
$ws.Range('D22').Value = '5.320'
$ws.Range('E22').Value = '  -2.74%  '
$ws.Range('E23').Value = '  -1.63%  '
$ws.Range('B24').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C24').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D24').Value = '2.132.19'
$ws.Range('E24').Value = '  -1.61%  '
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').Value = '2.092'
$ws.Range('E25').Value = '  -0.30%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').Value = '158.54'
$ws.Range('E26').Value = '  +0.88%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = '19.01'
$ws.Range('E27').Value = '  -1.83%  '
$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D28').Value = '5.660'
$ws.Range('E28').Value = '  -1.84%  '
$ws.Range('B29').Value = 'BitcoinCash'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D29').Value = '117.37'
$ws.Range('E29').Value = '  -0.91%  '
$ws.Range('B30').Value = 'LidoDAOToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D30').Value = '1.851'
$ws.Range('E30').Value = '  -0.56%  '
$ws.Range('B31').Value = 'Stellar'
$ws.Range('C31').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D31').Value = '0.09295'
$ws.Range('E31').Value = '  -0.71%  '
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').Value = '0.8685'
$ws.Range('E32').Value = '  +0.37%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = '5.079'
$ws.Range('E33').Value = '  -2.01%  '
$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D34').Value = '1.248'
$ws.Range('E34').Value = '  -4.39%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').Value = '3.030'
$ws.Range('E35').Value = '  -3.07%  '
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').Value = '0.05724'
$ws.Range('E36').Value = '  -0.76%  '
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').Value = '1.156'
$ws.Range('E37').Value = '  -0.09%  '
$ws.Range('B38').Value = 'Frax'
$ws.Range('C38').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D38').Value = '1.002'
$ws.Range('E38').Value = '  +0.08%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = '0.02041'
$ws.Range('E39').Value = '  -2.14%  '
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').Value = '0.5497'
$ws.Range('E40').Value = '  -2.82%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').Value = '7.395'
$ws.Range('E41').Value = '  -3.85%  '
$ws.Range('B42').Value = 'Algorand'
$ws.Range('C42').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D42').Value = '0.1751'
$ws.Range('E42').Value = '  -1.79%  '
$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D43').Value = '2.834'
$ws.Range('E43').Value = '  +3.68%  '
$ws.Range('B44').Value = 'Aptos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D44').Value = '9.310'
$ws.Range('E44').Value = '  -1.53%  '
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D45').Value = '0.5173'
$ws.Range('E45').Value = '  -2.21%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').Value = '11.23'
$ws.Range('E46').Value = '  -3.04%  '
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').Value = '0.06863'
$ws.Range('E47').Value = '  -0.12%  '
$ws.Range('D48').Value = '2.062'
$ws.Range('E48').Value = '  -2.14%  '
$ws.Range('B49').Value = 'PEPE'
$ws.Range('C49').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D49').Value = '0.000002605'
$ws.Range('E49').Value = '  -9.17%  '
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').Value = '1.781'
$ws.Range('E50').Value = '  -1.87%  '
$ws.Range('B51').Value = 'Quant'
$ws.Range('C51').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D51').Value = '110.54'
$ws.Range('E51').Value = '  -0.91%  '

$priceVolRange.Style = "Normal"

Write-Host "cryptos list updated"
